$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Swap the match-detail columns (F:V) between each of these row pairs. The
# identifying columns A (index) and E (date) stay put; only home/away teams,
# scores, odds, timestamps and url swap places between the two rows.
# ---------------------------------------------------------------------------
$swapPairs = @(
    @(60, 61),
    @(102, 103),
    @(106, 107),
    @(136, 137)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $rangeA = "F" + $r1 + ":V" + $r1
    $rangeB = "F" + $r2 + ":V" + $r2
    $valsA = $ws.Range($rangeA).Value()
    $valsB = $ws.Range($rangeB).Value()
    $ws.Range($rangeA).Value = $valsB
    $ws.Range($rangeB).Value = $valsA
}

# ---------------------------------------------------------------------------
# Append two new match rows (156, 157) at the bottom of the sheet.
# First clone formatting + the constant text columns (B,C,D) from the last
# existing row (155) via PasteSpecial so the bold/bordered index style (A),
# the datetime style (E) and the plain-text typing of B/C/D are preserved
# exactly, then fill in the per-row values.
# ---------------------------------------------------------------------------
$ws.Range("A155:E155").Copy()
$ws.Range("A156").PasteSpecial(-4122)
$ws.Range("A155:E155").Copy()
$ws.Range("A157").PasteSpecial(-4122)

$ws.Range("B155:D155").Copy()
$ws.Range("B156").PasteSpecial(-4104)
$ws.Range("B155:D155").Copy()
$ws.Range("B157").PasteSpecial(-4104)

$excel.CutCopyMode = $false

$ws.Range("A156").Value = 155
$ws.Range("E156").Value = 45233.95833333334
$ws.Range("F156").Value = "Arsenal Sarandi"
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = "Gimnasia L.P."
$ws.Range("I156").Value = 0
$ws.Range("J156").Value = 2.79
$ws.Range("K156").Value = "30/10/2023 02:12"
$ws.Range("L156").Value = 3.52
$ws.Range("M156").Value = "03/11/2023 22:58"
$ws.Range("N156").Value = 2.97
$ws.Range("O156").Value = "30/10/2023 02:12"
$ws.Range("P156").Value = 3.09
$ws.Range("Q156").Value = "03/11/2023 22:54"
$ws.Range("R156").Value = 2.87
$ws.Range("S156").Value = "30/10/2023 02:12"
$ws.Range("T156").Value = 2.35
$ws.Range("U156").Value = "03/11/2023 22:58"
$ws.Range("V156").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/arsenal-sarandi-gimnasia-l-p/dWerM2Lk/"

$ws.Range("A157").Value = 156
$ws.Range("E157").Value = 45234.04166666666
$ws.Range("F157").Value = "River Plate"
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = "Huracan"
$ws.Range("I157").Value = 2
$ws.Range("J157").Value = 1.5
$ws.Range("K157").Value = "31/10/2023 01:13"
$ws.Range("L157").Value = 1.48
$ws.Range("M157").Value = "04/11/2023 00:42"
$ws.Range("N157").Value = 4.21
$ws.Range("O157").Value = "31/10/2023 01:13"
$ws.Range("P157").Value = 4.3
$ws.Range("Q157").Value = "04/11/2023 00:50"
$ws.Range("R157").Value = 7.08
$ws.Range("S157").Value = "31/10/2023 01:13"
$ws.Range("T157").Value = 7.83
$ws.Range("U157").Value = "04/11/2023 00:50"
$ws.Range("V157").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/river-plate-huracan/4denLMzd/"
